$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3962581.2
$ws.Range("I40").Value = 50999.75
$ws.Range("K40").Value = 50999.75
$ws.Range("M40").Value = -50824.75
$ws.Range("H43").Value = 5131205
$ws.Range("J43").Value = 4500
$ws.Range("L43").Value = 4500
$ws.Range("N43").Value = -4638
$ws.Range("H92").Value = 606.5454999999999
$ws.Range("I92").Value = 606.5454999999999
$ws.Range("K92").Value = 606.5454999999999
$ws.Range("M92").Value = 641.4545000000001
$ws.Range("H96").Value = 3414.4546
$ws.Range("I96").Value = 470.88235
$ws.Range("J96").Value = 13422.6
$ws.Range("K96").Value = 1412.64705
$ws.Range("L96").Value = 40267.8
$ws.Range("M96").Value = -39.64705000000004
$ws.Range("N96").Value = -43013.8
$ws.Range("H103").Value = 541.1667
$ws.Range("I103").Value = 541.1667
$ws.Range("K103").Value = 1623.5001
$ws.Range("M103").Value = -1037.5001
$ws.Range("H106").Value = 25643432
$ws.Range("I106").Value = 27779550
$ws.Range("K106").Value = 27779550
$ws.Range("M106").Value = -27778919
$ws.Range("H116").Value = 19349948
$ws.Range("I116").Value = 28513162
$ws.Range("J116").Value = 5388.222
$ws.Range("K116").Value = 28513162
$ws.Range("L116").Value = 5388.222
$ws.Range("M116").Value = -28509720
$ws.Range("N116").Value = -12272.222
$ws.Range("H132").Value = 153735.16
$ws.Range("I132").Value = 224658.62
$ws.Range("J132").Value = 18642.857
$ws.Range("K132").Value = 673975.86
$ws.Range("L132").Value = 55928.571
$ws.Range("M132").Value = -671445.86
$ws.Range("N132").Value = -60988.571
$ws.Range("H135").Value = 3507.9778
$ws.Range("I135").Value = 946.5357
$ws.Range("J135").Value = 7726.8237
$ws.Range("K135").Value = 8518.8213
$ws.Range("L135").Value = 69541.4133
$ws.Range("M135").Value = -5983.8213
$ws.Range("N135").Value = -74611.4133

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14252.691
$ws.Range("I32").Value = 13541.491
$ws.Range("K32").Value = 13541.491
$ws.Range("M32").Value = -13254.491
$ws.Range("H61").Value = 22037.223
$ws.Range("I61").Value = 33267
$ws.Range("K61").Value = 33267
$ws.Range("M61").Value = -33055
$ws.Range("H63").Value = 4873.5
$ws.Range("I63").Value = 3250
$ws.Range("J63").Value = 6497
$ws.Range("K63").Value = 3250
$ws.Range("L63").Value = 6497
$ws.Range("M63").Value = -2564
$ws.Range("N63").Value = -7869
$ws.Range("H66").Value = 4873.5
$ws.Range("I66").Value = 3250
$ws.Range("J66").Value = 6497
$ws.Range("K66").Value = 16250
$ws.Range("L66").Value = 32485
$ws.Range("M66").Value = -12818
$ws.Range("N66").Value = -39349
$ws.Range("H132").Value = 18424.564
$ws.Range("I132").Value = 32218.842
$ws.Range("K132").Value = 96656.526
$ws.Range("M132").Value = -94126.526
$ws.Range("H136").Value = 22037.223
$ws.Range("I136").Value = 33267
$ws.Range("K136").Value = 99801
$ws.Range("M136").Value = -97251

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4389.087
$ws.Range("I20").Value = 4180
$ws.Range("K20").Value = 4180
$ws.Range("M20").Value = -3933
$ws.Range("H22").Value = 786.1111
$ws.Range("I22").Value = 813.7143
$ws.Range("K22").Value = 813.7143
$ws.Range("M22").Value = -640.7143
$ws.Range("H86").Value = 3568.8
$ws.Range("I86").Value = 2183.2
$ws.Range("J86").Value = 4954.4
$ws.Range("K86").Value = 2183.2
$ws.Range("L86").Value = 4954.4
$ws.Range("M86").Value = -1060.2
$ws.Range("N86").Value = -7200.4
$ws.Range("H89").Value = 3568.8
$ws.Range("I89").Value = 2183.2
$ws.Range("J89").Value = 4954.4
$ws.Range("K89").Value = 10916
$ws.Range("L89").Value = 24772
$ws.Range("M89").Value = -5300
$ws.Range("N89").Value = -36004
$ws.Range("H94").Value = 2741207.5
$ws.Range("I94").Value = 2741207.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2741207.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2740756.5
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 1626.5834
$ws.Range("I105").Value = 1626.5834
$ws.Range("K105").Value = 1626.5834
$ws.Range("M105").Value = 120.4166
$ws.Range("H132").Value = 120419
$ws.Range("J132").Value = 120419
$ws.Range("L132").Value = 120419
$ws.Range("N132").Value = -130539
$ws.Range("H134").Value = 3175.1667
$ws.Range("I134").Value = 2762.75
$ws.Range("K134").Value = 8288.25
$ws.Range("M134").Value = -5753.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2302.25
$ws.Range("I16").Value = 2302.25
$ws.Range("K16").Value = 2302.25
$ws.Range("M16").Value = -2015.25
$ws.Range("H22").Value = 515.1429000000001
$ws.Range("I22").Value = 241.46666
$ws.Range("K22").Value = 241.46666
$ws.Range("M22").Value = 108.53334
$ws.Range("H31").Value = 4578.7676
$ws.Range("I31").Value = 1618.2
$ws.Range("J31").Value = 7153.174
$ws.Range("K31").Value = 1618.2
$ws.Range("L31").Value = 7153.174
$ws.Range("M31").Value = -1323.2
$ws.Range("N31").Value = -7743.174
$ws.Range("H34").Value = 4578.7676
$ws.Range("I34").Value = 1618.2
$ws.Range("J34").Value = 7153.174
$ws.Range("K34").Value = 1618.2
$ws.Range("L34").Value = 7153.174
$ws.Range("M34").Value = -1416.2
$ws.Range("N34").Value = -7557.174
$ws.Range("H58").Value = 2642.3635
$ws.Range("I58").Value = 1449.9584
$ws.Range("J58").Value = 5822.1113
$ws.Range("K58").Value = 1449.9584
$ws.Range("L58").Value = 5822.1113
$ws.Range("M58").Value = -1246.9584
$ws.Range("N58").Value = -6228.1113
$ws.Range("H113").Value = 2302.25
$ws.Range("I113").Value = 2302.25
$ws.Range("K113").Value = 2302.25
$ws.Range("M113").Value = -132.25
$ws.Range("H136").Value = 2642.3635
$ws.Range("I136").Value = 1449.9584
$ws.Range("J136").Value = 5822.1113
$ws.Range("K136").Value = 4349.8752
$ws.Range("L136").Value = 17466.3339
$ws.Range("M136").Value = -1799.8752
$ws.Range("N136").Value = -22566.3339
$ws.Range("H141").Value = 79786.57000000001
$ws.Range("J141").Value = 90389
$ws.Range("L141").Value = 90389
$ws.Range("N141").Value = -100749

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1401.7778
$ws.Range("J22").Value = 1401.7778
$ws.Range("L22").Value = 4205.3334
$ws.Range("N22").Value = -4543.3334
$ws.Range("H27").Value = 1401.7778
$ws.Range("J27").Value = 1401.7778
$ws.Range("L27").Value = 4205.3334
$ws.Range("N27").Value = -4409.3334
$ws.Range("H34").Value = 939964.25
$ws.Range("I34").Value = 1291575.9
$ws.Range("K34").Value = 3874727.7
$ws.Range("M34").Value = -3874643.7
$ws.Range("H44").Value = 20000580
$ws.Range("I44").Value = 33333666
$ws.Range("K44").Value = 100000998
$ws.Range("M44").Value = -100000600
$ws.Range("H46").Value = 4994.6665
$ws.Range("I46").Value = 4994
$ws.Range("K46").Value = 14982
$ws.Range("M46").Value = -14891
$ws.Range("H124").Value = 354.33334
$ws.Range("I124").Value = 354.33334
$ws.Range("K124").Value = 1063.00002
$ws.Range("M124").Value = 3846.99998
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H129").Value = 1125.5
$ws.Range("J129").Value = 1997.5
$ws.Range("L129").Value = 5992.5
$ws.Range("N129").Value = -15992.5
$ws.Range("H131").Value = 31872596
$ws.Range("I131").Value = 41667000
$ws.Range("J131").Value = 30305490
$ws.Range("K131").Value = 125001000
$ws.Range("L131").Value = 90916470
$ws.Range("M131").Value = -124995960
$ws.Range("N131").Value = -90926550

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 371316.03
$ws.Range("I122").Value = 649518.1
$ws.Range("J122").Value = 7513.3076
$ws.Range("K122").Value = 1948554.3
$ws.Range("L122").Value = 22539.9228
$ws.Range("M122").Value = -1946104.3
$ws.Range("N122").Value = -27439.9228
$ws.Range("H126").Value = 4642.7646
$ws.Range("I126").Value = 2303.0833
$ws.Range("K126").Value = 6909.249899999999
$ws.Range("M126").Value = -4439.249899999999
$ws.Range("H134").Value = 5000000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2652.625
$ws.Range("I16").Value = 1927.375
$ws.Range("K16").Value = 1927.375
$ws.Range("M16").Value = -1757.375
$ws.Range("H40").Value = 4798.7
$ws.Range("I40").Value = 4443
$ws.Range("K40").Value = 4443
$ws.Range("M40").Value = -4307
$ws.Range("H68").Value = 2276758.5
$ws.Range("I68").Value = 11364236
$ws.Range("J68").Value = 4889.125
$ws.Range("K68").Value = 11364236
$ws.Range("L68").Value = 4889.125
$ws.Range("M68").Value = -11363487
$ws.Range("N68").Value = -6387.125
$ws.Range("H71").Value = 2276758.5
$ws.Range("I71").Value = 11364236
$ws.Range("J71").Value = 4889.125
$ws.Range("K71").Value = 56821180
$ws.Range("L71").Value = 24445.625
$ws.Range("M71").Value = -56817436
$ws.Range("N71").Value = -31933.625
$ws.Range("H93").Value = 4124.25
$ws.Range("I93").Value = 4165.6665
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 4165.6665
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -2917.6665
$ws.Range("N93").Value = -6496
$ws.Range("H136").Value = 3647.7827
$ws.Range("I136").Value = 2069.6155
$ws.Range("J136").Value = 5699.4
$ws.Range("K136").Value = 6208.8465
$ws.Range("L136").Value = 17098.2
$ws.Range("M136").Value = -3658.8465
$ws.Range("N136").Value = -22198.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
